$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.232270002365112
$ws.Range("B1").Value = 2.501514673233032
$ws.Range("C1").Value = 4.268093109130859
$ws.Range("D1").Value = 2.591438055038452
$ws.Range("E1").Value = 1.080204367637634
